$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "[1, 0, 1, 0, 1, 0, 0]"
$ws.Range("E11").Value = "['Normal', 'HardwareFault', 'RegulationViolation']"

$ws.Range("D38").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E38").Value = "['HardwareFault']"

$ws.Range("D53").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E53").Value = "['Normal', 'HardwareFault']"

$ws.Range("D56").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "['Normal']"

$ws.Range("D58").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E58").Value = "['Normal', 'ParamViolation']"

$ws.Range("D61").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E61").Value = "['SoftwareFault']"

$ws.Range("D69").Value = "[1, 1, 0, 0, 0, 1, 0]"
$ws.Range("E69").Value = "['Normal', 'SurroundingEnvironment', 'CommunicationIssue']"

$ws.Range("D70").Value = "[1, 1, 0, 0, 0, 1, 0]"
$ws.Range("E70").Value = "['Normal', 'SurroundingEnvironment', 'CommunicationIssue']"

$ws.Range("D73").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal']"

$ws.Range("D82").Value = "[1, 1, 1, 0, 0, 0, 0]"
$ws.Range("E82").Value = "['Normal', 'SurroundingEnvironment', 'HardwareFault']"

$ws.Range("D88").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E88").Value = "['Normal']"

$ws.Range("D92").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E92").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"
